$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the completed latitude/longitude survey data (rows 2-26) ---
$ws.Range("B2").Value = 51.339237199999999
$ws.Range("C2").Value = 12.3757456
$ws.Range("B3").Value = 51.338879599999999
$ws.Range("C3").Value = 12.3744023
$ws.Range("B4").Value = 51.328572600000001
$ws.Range("C4").Value = 12.331050899999999
$ws.Range("B5").Value = 51.330156500000001
$ws.Range("C5").Value = 12.337664800000001
$ws.Range("B6").Value = 51.3332081
$ws.Range("C6").Value = 12.3392631
$ws.Range("B7").Value = 51.332988299999997
$ws.Range("C7").Value = 12.337859999999999
$ws.Range("B8").Value = 51.332929900000003
$ws.Range("C8").Value = 12.3397963
$ws.Range("B9").Value = 51.317073200000003
$ws.Range("C9").Value = 12.3754904
$ws.Range("B10").Value = 51.327552400000002
$ws.Range("C10").Value = 12.373935700000001
$ws.Range("B11").Value = 51.332919199999999
$ws.Range("C11").Value = 12.373799099999999
$ws.Range("B12").Value = 51.339290499999997
$ws.Range("C12").Value = 12.3581936
$ws.Range("B13").Value = 51.339855
$ws.Range("C13").Value = 12.3674306
$ws.Range("B14").Value = 51.340958200000003
$ws.Range("C14").Value = 12.373424
$ws.Range("B15").Value = 51.344570400000002
$ws.Range("C15").Value = 12.374700900000001
$ws.Range("B16").Value = 51.3451278
$ws.Range("C16").Value = 12.3709808
$ws.Range("B17").Value = 51.344601400000002
$ws.Range("C17").Value = 12.3710349
$ws.Range("B18").Value = 51.3404673
$ws.Range("C18").Value = 12.379617700000001
$ws.Range("B19").Value = 51.337847600000003
$ws.Range("C19").Value = 12.3814378
$ws.Range("B20").Value = 51.345033000000001
$ws.Range("C20").Value = 12.3906513
$ws.Range("B21").Value = 51.3390968
$ws.Range("C21").Value = 12.393767499999999
$ws.Range("B22").Value = 51.344477099999999
$ws.Range("C22").Value = 12.366763000000001
$ws.Range("B23").Value = 51.331298599999997
$ws.Range("C23").Value = 12.374067699999999
$ws.Range("B24").Value = 51.333154299999997
$ws.Range("C24").Value = 12.3622513
$ws.Range("B25").Value = 51.325598300000003
$ws.Range("C25").Value = 12.371049899999999
$ws.Range("B26").Value = 51.323686600000002
$ws.Range("C26").Value = 12.372990400000001

# --- Mark an explicit "No Fill" on the data columns now that entry is complete ---
$ws.Range("B2:B26").Interior.ColorIndex = -4142  # xlColorIndexNone
$ws.Range("C2").Interior.ColorIndex = -4142      # xlColorIndexNone
$ws.Range("C3:C26").Interior.ColorIndex = -4142  # xlColorIndexNone

# --- Reset view: scroll back to the left edge and move the selection to B4 ---
$ws.Activate()
$ws.Range("B4").Select()
